# Manual update: document the shield SDI fix.
#
# The "Smash DI" row's description paragraph currently ends with the run
# " first smash DI input.". Word records an appended edit like this as a
# brand-new run (its own <w:r>) placed right after the existing one --
# it does not fold the new sentence into the run it follows. We
# reconstruct that paragraph (every existing run byte-for-byte
# unchanged) with the new run appended just before the closing </w:p>,
# then push it back in with Range.InsertXML so the paragraph's run
# boundaries end up exactly as Word would leave them.

$d = $word.ActiveDocument

# Locate the paragraph that documents "Smash DI" (the only paragraph in
# the manual ending in "first smash DI input.").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*first smash DI input.*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the Smash DI description paragraph."
}

$expectedText = "A tilt input on the first frame of hitlag will not prevent smash DI from occurring on the second frame of hitlag. Additionally, the second frame after entering the smash DI range counts toward the first smash DI input."
$actualText = $target.Range.Text.TrimEnd([char]13, [char]7)
if ($actualText -ne $expectedText) {
    throw "Smash DI paragraph text did not match the expected content; aborting to avoid corrupting the document."
}

$newSentence = " This fix applies to shield smash DI."

$newRunXml = '<w:r><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve">' + $newSentence + '</w:t></w:r>'

# Exact current OOXML for the paragraph (every run unchanged), with the
# new run spliced in right before </w:p>.
$paraBodyXml = '<w:p w14:paraId="1D8F7FB6" w14:textId="6FAC9796" w:rsidR="008A6D6D" w:rsidRPr="000C62FF" w:rsidRDefault="008A6D6D" w:rsidP="008A6D6D"><w:pPr><w:spacing w:line="168" w:lineRule="auto"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr></w:pPr><w:r w:rsidRPr="00FE3215"><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>A tilt input on the first frame of hitlag will not prevent smash DI from occurring on the second frame of hitlag. Additionally, the second frame after entering the smash DI range count</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr="00FE3215"><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve"> toward </w:t></w:r><w:r w:rsidR="00653F2C"><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>the</w:t></w:r><w:r w:rsidRPr="00FE3215"><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t xml:space="preserve"> first smash DI input.</w:t></w:r>' + $newRunXml + '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)

Write-Output "Updated Smash DI description with shield SDI fix note."
